$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: add two new date columns (AB, AC) after the existing AA column ---
# Clone the formatting of the last existing header cell (AA1) onto the two new
# header cells, then overwrite their text.
$ws.Range("AA1").Copy($ws.Range("AB1"))
$ws.Range("AA1").Copy($ws.Range("AC1"))
$ws.Range("AB1").Value = "04-07_A"
$ws.Range("AC1").Value = "04-07_0"

# --- Data rows: duplicate the last "_A"/"_0" column pair (Z/AA) into the new
# pair (AB/AC), keeping whatever value+style each source cell had. ---
for ($r = 2; $r -le 170; $r++) {
    $zAddr  = "Z" + $r
    $aaAddr = "AA" + $r
    $abAddr = "AB" + $r
    $acAddr = "AC" + $r

    # AB gets a copy of Z (value + style)
    $ws.Range($zAddr).Copy($ws.Range($abAddr))
    # AC gets a copy of AA as it was before this edit (value + style)
    $ws.Range($aaAddr).Copy($ws.Range($acAddr))

    # The original AA cell: if it holds a numeric-looking piece of text,
    # it becomes a real number (matching the "_0" pattern used by every
    # other numeric day column); if it is blank it is left untouched.
    $aaVal = $ws.Range($aaAddr).Value()
    if ($aaVal -ne $null -and "$aaVal" -match '^-?[0-9]+(\.[0-9]+)?$') {
        $ws.Range($aaAddr).Value = [double]"$aaVal"
    }

    # Column A on this row is a numeric ID that, on this particular row,
    # had been stored as text; normalise it to a real number too.
    $aAddr = "A" + $r
    $aVal = $ws.Range($aAddr).Value()
    if ($aVal -ne $null -and "$aVal" -match '^-?[0-9]+(\.[0-9]+)?$') {
        $ws.Range($aAddr).Value = [double]"$aVal"
    }
}
